$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,20
$vals = @("ECs", "Fgf2", "Sdc2", "ECs", 3, 1, 1.005243333333333, 3.01573, 0.07224874268505826, 0.07224874268505825, 3, 1, 3.057109, 9.171327, 0.02694952608666365, 0.02694952608666365, 3.073138441523333, 27.65824597371, 0.001947069375719627, 0.001947069375719627)
for ($i = 0; $i -lt 20; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("A2:T2").Value = $arr

$arr = New-Object 'object[,]' 1,20
$vals = @("ECs", "Fgf2", "Sdc2", "FAPs", 3, 1, 1.005243333333333, 3.01573, 0.07224874268505826, 0.07224874268505825, 3, 1, 89.02755999999999, 267.08268, 0.7848102735793893, 0.7848102735793893, 89.49436117293332, 805.4492505563999, 0.05670155551242748, 0.05670155551242746)
for ($i = 0; $i -lt 20; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("A3:T3").Value = $arr

$arr = New-Object 'object[,]' 1,20
$vals = @("ECs", "Fgf2", "Sdc2", "M2", 3, 1, 1.005243333333333, 3.01573, 0.07224874268505826, 0.07224874268505825, 2, 0.6666666666666666, 0.184005, 0.5520149999999999, 0.0016220709001794, 0.0016220709001794, 0.18496979955, 1.66472819595, 0.0001171925830839823, 0.0001171925830839822)
for ($i = 0; $i -lt 20; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("A4:T4").Value = $arr

$arr = New-Object 'object[,]' 1,20
$vals = @("ECs", "Fgf2", "Sdc2", "sCs", 3, 1, 1.005243333333333, 3.01573, 0.07224874268505826, 0.07224874268505825, 3, 1, 21.16964733333333, 63.508942, 0.1866181294337677, 0.1866181294337677, 21.28064685085111, 191.52582165766, 0.01348292521382718, 0.01348292521382718)
for ($i = 0; $i -lt 20; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("A5:T5").Value = $arr

$arr = New-Object 'object[,]' 1,20
$vals = @("FAPs", "Fgf2", "Sdc2", "ECs", 3, 1, 10.25983933333333, 30.779518, 0.7373940889775011, 0.737394088977501, 3, 1, 3.057109, 9.171327, 0.02694952608666365, 0.02694952608666365, 31.36544716448734, 282.289024480386, 0.01987242123705074, 0.01987242123705074)
for ($i = 0; $i -lt 20; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("A6:T6").Value = $arr

$arr = New-Object 'object[,]' 1,20
$vals = @("FAPs", "Fgf2", "Sdc2", "FAPs", 3, 1, 10.25983933333333, 30.779518, 0.7373940889775011, 0.737394088977501, 3, 1, 89.02755999999999, 267.08268, 0.7848102735793893, 0.7848102735793893, 913.4084618386933, 8220.676156548239, 0.5787144567062572, 0.5787144567062571)
for ($i = 0; $i -lt 20; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("A7:T7").Value = $arr

$arr = New-Object 'object[,]' 1,20
$vals = @("FAPs", "Fgf2", "Sdc2", "M2", 3, 1, 10.25983933333333, 30.779518, 0.7373940889775011, 0.737394088977501, 2, 0.6666666666666666, 0.184005, 0.5520149999999999, 0.0016220709001794, 0.0016220709001794, 1.88786173653, 16.99075562877, 0.001196105493694704, 0.001196105493694703)
for ($i = 0; $i -lt 20; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("A8:T8").Value = $arr

$arr = New-Object 'object[,]' 1,20
$vals = @("FAPs", "Fgf2", "Sdc2", "sCs", 3, 1, 10.25983933333333, 30.779518, 0.7373940889775011, 0.737394088977501, 3, 1, 21.16964733333333, 63.508942, 0.1866181294337677, 0.1866181294337677, 217.1971803833285, 1954.774623449956, 0.1376111055404985, 0.1376111055404985)
for ($i = 0; $i -lt 20; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("A9:T9").Value = $arr

$arr = New-Object 'object[,]' 1,20
$vals = @("M2", "Fgf2", "Sdc2", "ECs", 3, 1, 0.383875, 1.151625, 0.02758982345723265, 0.02758982345723265, 3, 1, 3.057109, 9.171327, 0.02694952608666365, 0.02694952608666365, 1.173547717375, 10.561929456375, 0.0007435326669871361, 0.0007435326669871361)
for ($i = 0; $i -lt 20; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("A10:T10").Value = $arr

$arr = New-Object 'object[,]' 1,20
$vals = @("M2", "Fgf2", "Sdc2", "FAPs", 3, 1, 0.383875, 1.151625, 0.02758982345723265, 0.02758982345723265, 3, 1, 89.02755999999999, 267.08268, 0.7848102735793893, 0.7848102735793893, 34.175454595, 307.579091355, 0.02165277689547781, 0.02165277689547781)
for ($i = 0; $i -lt 20; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("A11:T11").Value = $arr

$arr = New-Object 'object[,]' 1,20
$vals = @("M2", "Fgf2", "Sdc2", "M2", 3, 1, 0.383875, 1.151625, 0.02758982345723265, 0.02758982345723265, 2, 0.6666666666666666, 0.184005, 0.5520149999999999, 0.0016220709001794, 0.0016220709001794, 0.07063491937499999, 0.6357142743749999, 0.00004475264977106409, 0.00004475264977106408)
for ($i = 0; $i -lt 20; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("A12:T12").Value = $arr

$arr = New-Object 'object[,]' 1,20
$vals = @("M2", "Fgf2", "Sdc2", "sCs", 3, 1, 0.383875, 1.151625, 0.02758982345723265, 0.02758982345723265, 3, 1, 21.16964733333333, 63.508942, 0.1866181294337677, 0.1866181294337677, 8.126498370083334, 73.13848533075, 0.005148761244996644, 0.005148761244996644)
for ($i = 0; $i -lt 20; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("A13:T13").Value = $arr

$arr = New-Object 'object[,]' 1,20
$vals = @("sCs", "Fgf2", "Sdc2", "ECs", 3, 1, 2.264687, 6.794061, 0.162767344880208, 0.162767344880208, 3, 1, 3.057109, 9.171327, 0.02694952608666365, 0.02694952608666365, 6.923395009882999, 62.310555088947, 0.004386502806906145, 0.004386502806906145)
for ($i = 0; $i -lt 20; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("A14:T14").Value = $arr

$arr = New-Object 'object[,]' 1,20
$vals = @("sCs", "Fgf2", "Sdc2", "FAPs", 3, 1, 2.264687, 6.794061, 0.162767344880208, 0.162767344880208, 3, 1, 89.02755999999999, 267.08268, 0.7848102735793893, 0.7848102735793893, 201.61955777372, 1814.57601996348, 0.1277414844652268, 0.1277414844652268)
for ($i = 0; $i -lt 20; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("A15:T15").Value = $arr

$arr = New-Object 'object[,]' 1,20
$vals = @("sCs", "Fgf2", "Sdc2", "M2", 3, 1, 2.264687, 6.794061, 0.162767344880208, 0.162767344880208, 2, 0.6666666666666666, 0.184005, 0.5520149999999999, 0.0016220709001794, 0.0016220709001794, 0.4167137314349999, 3.750423582914999, 0.0002640201736296498, 0.0002640201736296498)
for ($i = 0; $i -lt 20; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("A16:T16").Value = $arr

$arr = New-Object 'object[,]' 1,20
$vals = @("sCs", "Fgf2", "Sdc2", "sCs", 3, 1, 2.264687, 6.794061, 0.162767344880208, 0.162767344880208, 3, 1, 21.16964733333333, 63.508942, 0.1866181294337677, 0.1866181294337677, 47.94262511038467, 431.483625993462, 0.03037533743444537, 0.03037533743444536)
for ($i = 0; $i -lt 20; $i++) { $arr[0,$i] = $vals[$i] }
$ws.Range("A17:T17").Value = $arr

Write-Host "Updated rows 2-17 on sheet1"
